$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3774.0833
$ws.Range("I32").Value = 3923.5
$ws.Range("J32").Value = 3475.25
$ws.Range("K32").Value = 3923.5
$ws.Range("L32").Value = 3475.25
$ws.Range("M32").Value = -3597.5
$ws.Range("N32").Value = -4127.25
$ws.Range("H64").Value = 103440
$ws.Range("I64").Value = 4625
$ws.Range("J64").Value = 169316.67
$ws.Range("K64").Value = 4625
$ws.Range("L64").Value = 169316.67
$ws.Range("M64").Value = -4377
$ws.Range("N64").Value = -169812.67
$ws.Range("H67").Value = 103440
$ws.Range("I67").Value = 4625
$ws.Range("J67").Value = 169316.67
$ws.Range("K67").Value = 4625
$ws.Range("L67").Value = 169316.67
$ws.Range("M67").Value = -3767
$ws.Range("N67").Value = -171032.67
$ws.Range("H74").Value = 4357.143
$ws.Range("I74").Value = 4333.3335
$ws.Range("K74").Value = 4333.3335
$ws.Range("M74").Value = -3397.3335
$ws.Range("H76").Value = 3300
$ws.Range("I76").Value = 3283.3333
$ws.Range("K76").Value = 3283.3333
$ws.Range("M76").Value = -2968.3333
$ws.Range("H77").Value = 4357.143
$ws.Range("I77").Value = 4333.3335
$ws.Range("K77").Value = 21666.6675
$ws.Range("M77").Value = -16986.6675
$ws.Range("H79").Value = 3300
$ws.Range("I79").Value = 3283.3333
$ws.Range("K79").Value = 3283.3333
$ws.Range("M79").Value = -2191.3333
$ws.Range("H123").Value = 42632
$ws.Range("J123").Value = 42632
$ws.Range("L123").Value = 42632
$ws.Range("N123").Value = -52432
$ws.Range("H128").Value = 41860
$ws.Range("J128").Value = 41860
$ws.Range("L128").Value = 41860
$ws.Range("N128").Value = -51820
$ws.Range("H132").Value = 373942.72
$ws.Range("I132").Value = 388248.2
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 1164744.6
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1162214.6
$ws.Range("N132").Value = -11060

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1410.8572
$ws.Range("J64").Value = 763.1429000000001
$ws.Range("L64").Value = 763.1429000000001
$ws.Range("N64").Value = -1213.1429
$ws.Range("H67").Value = 1410.8572
$ws.Range("J67").Value = 763.1429000000001
$ws.Range("L67").Value = 763.1429000000001
$ws.Range("N67").Value = -2323.1429

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 38.07143
$ws.Range("I7").Value = 29.875
$ws.Range("J7").Value = 49
$ws.Range("K7").Value = 29.875
$ws.Range("L7").Value = 49
$ws.Range("M7").Value = 83.125
$ws.Range("N7").Value = -275
$ws.Range("H31").Value = 1371.7142
$ws.Range("I31").Value = 1240.3
$ws.Range("J31").Value = 4000
$ws.Range("K31").Value = 1240.3
$ws.Range("L31").Value = 4000
$ws.Range("M31").Value = -945.3
$ws.Range("N31").Value = -4590
$ws.Range("H34").Value = 1371.7142
$ws.Range("I34").Value = 1240.3
$ws.Range("J34").Value = 4000
$ws.Range("K34").Value = 1240.3
$ws.Range("L34").Value = 4000
$ws.Range("M34").Value = -1038.3
$ws.Range("N34").Value = -4404
$ws.Range("H41").Value = 7376.5
$ws.Range("I41").Value = 6491.8
$ws.Range("J41").Value = 11800
$ws.Range("K41").Value = 6491.8
$ws.Range("L41").Value = 11800
$ws.Range("M41").Value = -6063.8
$ws.Range("N41").Value = -12656
$ws.Range("H50").Value = 7794
$ws.Range("I50").Value = 7900
$ws.Range("J50").Value = 7758.6665
$ws.Range("K50").Value = 7900
$ws.Range("L50").Value = 7758.6665
$ws.Range("M50").Value = -7275
$ws.Range("N50").Value = -9008.666499999999
$ws.Range("H51").Value = 25500
$ws.Range("J51").Value = 25500
$ws.Range("L51").Value = 25500
$ws.Range("N51").Value = -26972
$ws.Range("H60").Value = 9613.111000000001
$ws.Range("I60").Value = 6633.3335
$ws.Range("J60").Value = 11103
$ws.Range("K60").Value = 6633.3335
$ws.Range("L60").Value = 11103
$ws.Range("M60").Value = -6122.3335
$ws.Range("N60").Value = -12125
$ws.Range("H61").Value = 25500
$ws.Range("J61").Value = 25500
$ws.Range("L61").Value = 25500
$ws.Range("N61").Value = -26196
$ws.Range("H107").Value = 2193278
$ws.Range("I107").Value = 3205404.5
$ws.Range("J107").Value = 337.66666
$ws.Range("K107").Value = 3205404.5
$ws.Range("L107").Value = 337.66666
$ws.Range("M107").Value = -3203484.5
$ws.Range("N107").Value = -4177.66666

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 731.3
$ws.Range("I98").Value = 900
$ws.Range("K98").Value = 2700
$ws.Range("M98").Value = -1202
$ws.Range("H122").Value = 27474608
$ws.Range("I122").Value = 66667256
$ws.Range("J122").Value = 2979203.5
$ws.Range("K122").Value = 600005304
$ws.Range("L122").Value = 26812831.5
$ws.Range("M122").Value = -600002854
$ws.Range("N122").Value = -26817731.5

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5289.3
$ws.Range("I70").Value = 5240.353
$ws.Range("J70").Value = 5566.6665
$ws.Range("K70").Value = 5240.353
$ws.Range("L70").Value = 5566.6665
$ws.Range("M70").Value = -4970.353
$ws.Range("N70").Value = -6106.6665
$ws.Range("H73").Value = 5289.3
$ws.Range("I73").Value = 5240.353
$ws.Range("J73").Value = 5566.6665
$ws.Range("K73").Value = 5240.353
$ws.Range("L73").Value = 5566.6665
$ws.Range("M73").Value = -4304.353
$ws.Range("N73").Value = -7438.6665
$ws.Range("H80").Value = 112090.09
$ws.Range("J80").Value = 153248.88
$ws.Range("L80").Value = 153248.88
$ws.Range("N80").Value = -155244.88
$ws.Range("H83").Value = 112090.09
$ws.Range("J83").Value = 153248.88
$ws.Range("L83").Value = 766244.4
$ws.Range("N83").Value = -776228.4
$ws.Range("H107").Value = 591.4
$ws.Range("J107").Value = 740.2857
$ws.Range("L107").Value = 740.2857
$ws.Range("N107").Value = -4580.2857
$ws.Range("H134").Value = 20161
$ws.Range("J134").Value = 20161
$ws.Range("L134").Value = 60483
$ws.Range("N134").Value = -65553

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H135").Value = 32666.666
$ws.Range("J135").Value = 32666.666
$ws.Range("L135").Value = 32666.666
$ws.Range("N135").Value = -42806.666
